$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Helper pattern used throughout:
#   - New header cells in column A are created by setting the text value
#     and then copying the (bold/border/centered) format that already
#     exists on the B1 header cell of the same row via PasteSpecial
#     (paste formats only), so the resulting cell reuses the existing
#     "header" cell style instead of creating a brand new style.
#   - Existing column-A label cells lose their header styling, so we use
#     ClearFormats() on them after (optionally) updating their text.
# -----------------------------------------------------------------------

# ------------------------------------------------------------------
# Sheets 1-4 share an identical table layout and an identical set of
# changes: add an "Fonte/Tecnologia" header in A1, strip the header
# style from A2:A12, and fix a few accented labels.
# ------------------------------------------------------------------
$sheetIndexes = 1,2,3,4

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Add new header cell A1, copying the existing header format from B1.
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)

    # Row labels: fix accents where needed and drop the header styling.
    $ws.Range("A2").ClearFormats()

    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A3").ClearFormats()

    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A4").ClearFormats()

    $ws.Range("A5").ClearFormats()

    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A6").ClearFormats()

    $ws.Range("A7").ClearFormats()

    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A8").ClearFormats()

    $ws.Range("A9").ClearFormats()

    $ws.Range("A10").ClearFormats()

    $ws.Range("A11").Value = "Pot. Compl."
    $ws.Range("A11").ClearFormats()

    $ws.Range("A12").ClearFormats()
}

# ------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
#   - add "Período" header in A1
#   - fix accents on A2/A3 and drop their header styling
#   - remove the "Teto" row entirely
# ------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A2").ClearFormats()

$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A3").ClearFormats()

$ws5.Rows.Item(4).Delete()

# ------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
#   - add "Tipo Expansão" header in A1
#   - change B1 from "Custo" to "2015" (stays a text header cell)
#   - fix accents on A2/A3 and drop their header styling
#   - update the B2/B3 values
# ------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# Preserve B1's existing header style in a scratch cell so we can set
# its value (forcing text, not a number) and then restore the style.
$ws6.Range("B1").Copy()
$ws6.Range("Z1").PasteSpecial(-4122)

$ws6.Range("B1").Value = "'2015"

$ws6.Range("Z1").Copy()
$ws6.Range("B1").PasteSpecial(-4122)
$ws6.Range("Z1").Clear()

# New header cell A1, reusing B1's (now restored) header style.
$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("A2").ClearFormats()
$ws6.Range("B2").Value = 471

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("A3").ClearFormats()
$ws6.Range("B3").Value = 99
